$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. "misrak meskan" (guraghe zone) is renamed to "meskan"
$ws.Range("C5").Value = "meskan"

# 2. The "missing" wereda row (guraghe zone) is removed entirely
$ws.Rows("6").Delete()

# 3. The whole "halaba" zone (atoti hullo, wera, wera dijjo) is removed
#    After step 2, these rows shifted up from 11:13 to 10:12
$ws.Rows("10:12").Delete()

# 4. A new wereda "alaba special" under a new "siltie" zone is inserted
#    right after the "silti" row (silite zone), which after the prior
#    deletions now sits at row 18
$ws.Rows("19").Insert()
$ws.Range("A19").Value = "snnp"
$ws.Range("B19").Value = "siltie"
$ws.Range("C19").Value = "alaba special"
